# v1.0.19 corrección Excel threshold under 10
# Update the Send_To / Send_To2 / Notify_To email list so it no longer
# includes lucy.serrano@vcimentos.com, keeping only jmruiz@rpatechnologies.es

$wb = $excel.ActiveWorkbook

$oldValue = "lucy.serrano@vcimentos.com;jmruiz@rpatechnologies.es"
$newValue = "jmruiz@rpatechnologies.es"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ($cell.Value2 -eq $oldValue) {
            $cell.Value = $newValue
        }
    }
}
